# Add a new "Compact List" paragraph style (styleId "CompactList"),
# mirroring the existing "Compact" style: based on Body Text, quick
# style, with 1.8pt (36 twips) spacing before/after.

$d = $word.ActiveDocument

$compactList = $d.Styles.Add("Compact List", 1)
$compactList.BaseStyle = "BodyText"
$compactList.QuickStyle = $true
$compactList.ParagraphFormat.SpaceBefore = 1.8
$compactList.ParagraphFormat.SpaceAfter = 1.8
